$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Insert a new row at position 36 (shifts existing rows 36+ down to 37+)
$ws.Rows.Item(36).Insert()

# Populate the new row 36 with the "r146a" rule text
$ws.Cells.Item(36, 1).Value = "r146a"
$ws.Cells.Item(36, 2).Value = "<Bold>e146a Stealing Count Drogat Jewels</Bold>`r`n<LineBreak/><LineBreak/>Using the foulbane in Drogat Castle, you can spend a day instead of a normal daily action in arranging for a special theft of the Count's personnel jewels.`r`n<LineBreak/><LineBreak/>At the end of the day, you escape from the hex`r`n <InlineUIContainer><Button Content='r218' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>`r`n with wealth 110. However, you can never return to the castle hex due to the Count's anger. In the process of the theft, you might acquire magic items as part of the wealth 110 check. `r`n<LineBreak/><LineBreak/>`r`n                      <InlineUIContainer><Image Source='../bin/Images/CountDrogatJewels.gif' Name='CountDrogatJewels' Height='300' Width='300'></Image></InlineUIContainer>"

# Copy cell style formatting from row above (row 35) to keep consistent styles
$ws.Cells.Item(35, 1).Copy()
$ws.Cells.Item(36, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(35, 2).Copy()
$ws.Cells.Item(36, 2).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Set the row height for the new row
$ws.Rows.Item(36).RowHeight = 105

# Update the view: scroll position and active selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 33
$ws.Range("B36").Select()
